$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 452 (shifts old rows 452:533 down to 456:537)
$insertRange = $ws.Range("A452:T455")
$insertRange.Insert()

# Populate the 4 newly inserted rows with data
$rows = @(452, 453, 454, 455)
$L = @('Especial', 'Extra (doble especial)', 'Primera', 'Segunda')
$M = @(200, 150, 160, 200)
$N = @(18000, 21600, 14400, 10800)
$S = @(1000, 1200, 800, 600)

for ($i = 0; $i -lt 4; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = 'Vega Central Mapocho de Santiago'
    $ws.Cells.Item($r, 3).Value = 'Metropolitana'
    $ws.Cells.Item($r, 4).Value = 44694
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 'Fruta'
    $ws.Cells.Item($r, 7).Value = 100101
    $ws.Cells.Item($r, 8).Value = 'Berries'
    $ws.Cells.Item($r, 9).Value = 100101007
    $ws.Cells.Item($r, 10).Value = 'Kiwi'
    $ws.Cells.Item($r, 11).Value = 'Hayward'
    $ws.Cells.Item($r, 12).Value = $L[$i]
    $ws.Cells.Item($r, 13).Value = $M[$i]
    $ws.Cells.Item($r, 14).Value = $N[$i]
    $ws.Cells.Item($r, 15).Value = $N[$i]
    $ws.Cells.Item($r, 16).Value = $N[$i]
    $ws.Cells.Item($r, 17).Value = '$/caja 18 kilos'
    $ws.Cells.Item($r, 18).Value = 'Provincia de Curicó'
    $ws.Cells.Item($r, 19).Value = $S[$i]
    $ws.Cells.Item($r, 20).Value = 18
}

Write-Host "Done"
